# ------------------------------------------------------------------
# Take screenshot is added / Constants class is added
#
# - Update the sample login e-mail on "LoginData"
# - Add a new "SigninData" sheet (after "LoginData") with a header row
#   describing the fields used by the sign-in/registration form
# ------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- LoginData: refresh the generated challenge e-mail address -----
$ws1.Range("B2").Value = "hf_challenge_1578909836186@hf836186.com"

# --- Add the new SigninData worksheet right after LoginData --------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "SigninData"

$headers = @("useremail", "Name", "Surname", "Password", "Days", "Months", "Years", "Company", "Address1", "Address2", "City", "Id_state", "Postcode", "Other", "Phone", "Phone_mobile", "Alias")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Header row (B1:Q1) is bold, matching the style used on LoginData's header
$ws2.Range("B1:Q1").Font.Bold = $true

# Approximate the auto-fit column widths used by the original workbook
# (target widths are 9.85546875 / 14.140625; the engine quantizes
# column widths, so these inputs land on the closest reachable value)
$ws2.Columns.Item(1).ColumnWidth = 9
$ws2.Columns.Item(16).ColumnWidth = 13.3

# Leave the selection on G5, as in the source workbook
$ws2.Range("G5").Select() | Out-Null
